$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.7376376588883126
$ws.Range("J2").Value = 0.7376376588883125
$ws.Range("M2").Value = 3.618510333333333
$ws.Range("N2").Value = 10.855531
$ws.Range("O2").Value = 0.1815566256530994
$ws.Range("P2").Value = 0.1815566256530994
$ws.Range("Q2").Value = 2.108052451271555
$ws.Range("R2").Value = 18.972472061444
$ws.Range("S2").Value = 0.133923004302414
$ws.Range("T2").Value = 0.133923004302414

# Row 3
$ws.Range("I3").Value = 0.7376376588883126
$ws.Range("J3").Value = 0.7376376588883125
$ws.Range("O3").Value = 0.1937079481987336
$ws.Range("P3").Value = 0.1937079481987336
$ws.Range("S3").Value = 0.1428862774173724
$ws.Range("T3").Value = 0.1428862774173724

# Row 4
$ws.Range("I4").Value = 0.7376376588883126
$ws.Range("J4").Value = 0.7376376588883125
$ws.Range("M4").Value = 12.45127566666667
$ws.Range("N4").Value = 37.353827
$ws.Range("O4").Value = 0.6247354261481669
$ws.Range("P4").Value = 0.6247354261481669
$ws.Range("Q4").Value = 7.253797771083111
$ws.Range("R4").Value = 65.284179939748
$ws.Range("S4").Value = 0.4608283771685261
$ws.Range("T4").Value = 0.460828377168526

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2072096666666667
$ws.Range("H5").Value = 0.621629
$ws.Range("I5").Value = 0.2623623411116874
$ws.Range("J5").Value = 0.2623623411116874
$ws.Range("M5").Value = 3.618510333333333
$ws.Range("N5").Value = 10.855531
$ws.Range("O5").Value = 0.1815566256530994
$ws.Range("P5").Value = 0.1815566256530994
$ws.Range("Q5").Value = 0.7497903199998888
$ws.Range("R5").Value = 6.748112879998999
$ws.Range("S5").Value = 0.0476336213506854
$ws.Range("T5").Value = 0.04763362135068541

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.2072096666666667
$ws.Range("H6").Value = 0.621629
$ws.Range("I6").Value = 0.2623623411116874
$ws.Range("J6").Value = 0.2623623411116874
$ws.Range("O6").Value = 0.1937079481987336
$ws.Range("P6").Value = 0.1937079481987336
$ws.Range("Q6").Value = 0.7999727024226667
$ws.Range("R6").Value = 7.199754321804
$ws.Range("S6").Value = 0.05082167078136123
$ws.Range("T6").Value = 0.05082167078136123

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.2072096666666667
$ws.Range("H7").Value = 0.621629
$ws.Range("I7").Value = 0.2623623411116874
$ws.Range("J7").Value = 0.2623623411116874
$ws.Range("M7").Value = 12.45127566666667
$ws.Range("N7").Value = 37.353827
$ws.Range("O7").Value = 0.6247354261481669
$ws.Range("P7").Value = 0.6247354261481669
$ws.Range("Q7").Value = 2.580024680464777
$ws.Range("R7").Value = 23.220222124183
$ws.Range("S7").Value = 0.1639070489796408
$ws.Range("T7").Value = 0.1639070489796408
